$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Row 23 currently carries the "highlighted last row" formatting
#    (B23/C23). That formatting needs to move down to the new last
#    row (25), while row 23 itself reverts to the normal row format
#    used by the rows above it (e.g. row 22).
# ------------------------------------------------------------------
$ws.Range("B22:C22").Copy()
$ws.Range("B23:C25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A22").Copy()
$ws.Range("A24:A25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Fill in the two new log entries.
# ------------------------------------------------------------------
$ws.Range("A24").Value = 45713
$ws.Range("B24").Value = 1.5
$ws.Range("C24").Value = "Researched about AI tool for remediation suggestion"

$ws.Range("A25").Value = 45714
$ws.Range("B25").Value = 4
$ws.Range("C25").Value = "Tried to integrate LangChain (RAG AI tool) with existing code"

# ------------------------------------------------------------------
# 3. Reflect the author's final selection state when the file was
#    last saved (cell E22 was the active cell).
# ------------------------------------------------------------------
$ws.Range("E22").Select()
